# "add say command, just chatting in game"
# Row 9  ("Добавить чат (say)"): progress moves from "В процессе" to "Cделано"
# Row 15 ("Сделать усложнение врагов на каждом новом этаже...") and
# Row 24 ("Сделать отображение других игроков на карте"): progress reverted
# from "В процессе" back to "Не сделано" (version cleared)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: match the "Cделано" look already used by B2 (fill + border), then
# set the cell's own text.
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("B9").Value = "Cделано"

# Rows 15 & 24: match the "Не сделано" / empty-version look already used by
# B5:C5, then set the cells' own text/empty value.
$ws.Range("B5:C5").Copy()

$ws.Range("B15:C15").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("B15").Value = "Не сделано"
$ws.Range("C15").Value = ""

$ws.Range("B5:C5").Copy()
$ws.Range("B24:C24").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("B24").Value = "Не сделано"
$ws.Range("C24").Value = ""

# View changes: zoom level and active selection
$excel.ActiveWindow.Zoom = 130
$ws.Range("E8").Select()
